# ToDo.docx edit:
#  - mark several completed to-do items with strikethrough
#  - re-color two "homepage" bullet items from black/text1 to dark red (C00000)
#  - relocate the Word "_GoBack" last-edit bookmark (best effort; some
#    headless COM hosts treat _GoBack as an immutable/auto-managed mark)

$d = $word.ActiveDocument

function Get-ParagraphByText($needle) {
    $rng = $d.Content
    $null = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    return $rng.Paragraphs(1)
}

function Set-ParagraphStrike($needle) {
    $para = Get-ParagraphByText $needle
    $para.Range.Font.StrikeThrough = 1
}

function Set-ParagraphColor($needle, $wdColor) {
    $para = Get-ParagraphByText $needle
    $para.Range.Font.Color = $wdColor
}

# 1) "saitis meilis sistemashi..." (email-notification intro bullet) - done -> strike
Set-ParagraphStrike "საიტის მეილის სისტემაში უდნა იყოს 3 ტიპის შეტყობინება"

# 2) "Challenge - მეგობარი გამოიწვევს..." bullet - done -> strike
Set-ParagraphStrike "Challenge - მეგობარი გამოიწვევს მეგობარს"

# 3) "მომხმარებლის მეგობრების აქტივობა..." bullet - done -> strike
Set-ParagraphStrike "მომხმარებლის მეგობრების აქტივობა"

# 4) "წინა დღის განმავლობაში..." bullet - recolor black -> C00000 (192)
Set-ParagraphColor "წინა დღის განმავლობაში საუკეთესო შედეგების მქონეთა სია" 192

# 5) "შემაჯამებელი სტატისტიკა;" bullet - recolor black -> C00000 (192)
Set-ParagraphColor "შემაჯამებელი სტატისტიკა" 192

# 6) "თუ მომხმარებელი შედის საწყის (login) გვერდზე..." bullet - done -> strike
Set-ParagraphStrike "თუ მომხმარებელი შედის საწყის"

# 7) "რამდენიმე მომხმარებელს ერთად უნდა შეეძოს საიტზე შესვლა..." bullet - done -> strike
Set-ParagraphStrike "რამდენიმე მომხმარებელს ერთად უნდა შეეძოს საიტზე შესვლა"

# 8) Move the "_GoBack" bookmark (Word's last-edit-location marker) from the
#    "დაურეგისტრირებელ მომხმარებელსაც..." bullet to the "წინა დღის..." bullet,
#    mirroring where the author's cursor was when the file was last saved.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$target = Get-ParagraphByText "წინა დღის განმავლობაში საუკეთესო შედეგების მქონეთა სია"
$bmRange = $d.Range($target.Range.Start, $target.Range.End)
$d.Bookmarks.Add("_GoBack", $bmRange)
